$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.899.15"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "2.826.04"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'502.96"
$ws.Range("E5").Value = "  -5.24%  "
$ws.Range("D6").Value = "'134.63"
$ws.Range("E6").Value = "  -7.65%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  -6.24%  "
$ws.Range("D9").Value = "2.823.99"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -6.44%  "
$ws.Range("D11").Value = "'5.90"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'0.345"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "3.320.96"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").Value = "59.027.65"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("E16").Value = "  -7.13%  "
$ws.Range("D17").Value = "2.827.02"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  -5.88%  "
$ws.Range("D19").Value = "'4.67"
$ws.Range("E19").Value = "  -6.74%  "
$ws.Range("D20").Value = "'10.94"
$ws.Range("E20").Value = "  -6.48%  "
$ws.Range("D21").Value = "'346.65"
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("D22").Value = "'6.21"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'62.87"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("D26").Value = "'0.424"
$ws.Range("E26").Value = "  -6.92%  "
$ws.Range("D27").Value = "'0.169"
$ws.Range("E27").Value = "  -8.30%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'7.30"
$ws.Range("E29").Value = "  -6.10%  "
$ws.Range("D30").Value = "0.0₃0791"
$ws.Range("E30").Value = "  -10.08%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").Value = "'18.97"
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("D34").Value = "'150.85"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").Value = "'5.28"
$ws.Range("E36").Value = "  -6.80%  "
$ws.Range("D37").Value = "'0.897"
$ws.Range("E37").Value = "  -11.42%  "
$ws.Range("D38").Value = "'1.11"
$ws.Range("E38").Value = "  -9.13%  "
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "2.229.44"
$ws.Range("E40").Value = "  -5.63%  "
$ws.Range("D41").Value = "'0.625"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Value = "'3.49"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "'0.0553"
$ws.Range("E44").Value = "  -3.80%  "
$ws.Range("E45").Value = "  -10.93%  "
$ws.Range("D46").Value = "'19.08"
$ws.Range("E46").Value = "  -9.27%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.56"
$ws.Range("E49").Value = "  -7.90%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0885"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").Value = "'17.29"
$ws.Range("E51").Value = "  -7.98%  "
